$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted ahead of what was row 152, pushing
# every subsequent record down by one row (old row 152 -> new row 153, ...,
# old row 231 -> new row 232).
$ws.Rows.Item(152).Insert()

# Populate the newly inserted row 152 with the new record's data.
$ws.Range("A152").Value = 8
$ws.Range("B152").Value = "Terminal La Palmera de La Serena"
$ws.Range("C152").Value = "Coquimbo"
$ws.Range("D152").Value = 44460
$ws.Range("E152").Value = 4
$ws.Range("F152").Value = 100114001
$ws.Range("G152").Value = "Papa"
$ws.Range("H152").Value = "Cardinal"
$ws.Range("I152").Value = "1a (cosecha)"
$ws.Range("J152").Value = 2600
$ws.Range("K152").Value = 12000
$ws.Range("L152").Value = 12500
$ws.Range("M152").Value = 12250
$ws.Range("N152").Value = "`$/saco 25 kilos"
$ws.Range("O152").Value = "Provincia del Elquí"
$ws.Range("P152").Value = 490
$ws.Range("Q152").Value = 25
$ws.Range("R152").Value = "Hortaliza"
